$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update existing value in A2 (was 259 -> now 388)
$ws.Range("A2").Value = 388

# Shift/replace values for rows 6-9 and add a new row 10
$ws.Range("A6").Value = 780
$ws.Range("A7").Value = 2200
$ws.Range("A8").Value = 2201
$ws.Range("A9").Value = 2447
$ws.Range("A10").Value = 119
